$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values scraped from coinranking.com.
# Price cells are forced to text format ("@") since values like "51.735.98"
# or "0.933" must remain literal text and not be reinterpreted as numbers/dates.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.735.98'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.782.06'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '356.51'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.44'
$ws.Range("E6").Value = '  -2.52%  '
$ws.Range("E7").Value = '  -2.94%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("E11").Value = '  +3.61%  '
$ws.Range("E12").Value = '  -2.11%  '
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.59'
$ws.Range("E14").Value = '  -2.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.219.45'
$ws.Range("E15").Value = '  -1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.781.46'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.933'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.689.14'
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.09'
$ws.Range("E20").Value = '  -3.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.21'
$ws.Range("E21").Value = '  -1.86%  '
$ws.Range("E22").Value = '  -2.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.19'
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.46'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("E25").Value = '  -2.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.36'
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.167'
$ws.Range("E28").Value = '  +15.59%  '
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.17'
$ws.Range("E30").Value = '  -4.12%  '
$ws.Range("E31").Value = '  +3.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.75'
$ws.Range("E32").Value = '  -0.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.68'
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0446'
$ws.Range("E34").Value = '  -7.96%  '
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.19'
$ws.Range("E36").Value = '  -7.19%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.89'
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("E39").Value = '  -4.61%  '
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("E43").Value = '  -2.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.88'
$ws.Range("E44").Value = '  -6.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.66'
$ws.Range("E45").Value = '  -6.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.091.66'
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.950'
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("E50").Value = '  -7.04%  '
$ws.Range("E51").Value = '  -7.10%  '
